# Adds a new "2022-Q4" sheet (with fund holding detail) right after the
# "总计" (summary) sheet, and inserts a corresponding summary row at the
# top of the "总计" sheet's data table.

$wb = $excel.ActiveWorkbook

# Helper: write a value into a cell while forcing it to be stored as TEXT
# (shared/inline string) even when it looks like a number (e.g. "1.80",
# "0.00", "016470"). Excel's normal smart-type detection would otherwise
# silently convert such strings to numbers and drop the formatting
# (leading/trailing zeros). We flip the cell to the "Text" number format
# before assigning the value, then clear the format again so no stray
# style index is left behind on the cell.
function Set-TextValue($cell, [string]$val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# Helper: write a genuine numeric value into a cell.
function Set-NumberValue($cell, $val) {
    $cell.Value = $val
}

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet right after "总计" and before
#    "2022-Q3".
# ---------------------------------------------------------------------
$summarySheet = $wb.Worksheets.Item("总计")
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheetName = $lastSheet.Name

$newSheet = $wb.Worksheets.Add([Type]::Missing, $summarySheet)
$newSheet.Name = "2022-Q4"

# Header row
Set-TextValue $newSheet.Cells.Item(1,2) "基金代码"
Set-TextValue $newSheet.Cells.Item(1,3) "基金名称"
Set-TextValue $newSheet.Cells.Item(1,4) "基金规模"
Set-TextValue $newSheet.Cells.Item(1,5) "股票总仓位"
Set-TextValue $newSheet.Cells.Item(1,6) "仓位占比"
Set-TextValue $newSheet.Cells.Item(1,7) "持有市值(亿元)"
Set-TextValue $newSheet.Cells.Item(1,8) "仓位排名"

$rows = @(
    @{ idx = 0; code = "001092"; name = "广发纳斯达克生物科技指数人民币A";            scale = "1.80"; pos = "90.20"; pct = "7.16"; mv = "0.1289"; mvIsZero = $false; rank = 1 },
    @{ idx = 1; code = "001093"; name = "广发纳斯达克生物科技指数美元A";            scale = "1.80"; pos = "90.20"; pct = "7.16"; mv = "0.1289"; mvIsZero = $false; rank = 1 },
    @{ idx = 2; code = "513290"; name = "汇添富纳斯达克生物科技ETF（QDII）";        scale = "0.75"; pos = "99.38"; pct = "8.01"; mv = "0.0601"; mvIsZero = $false; rank = 1 },
    @{ idx = 3; code = "016470"; name = "广发纳斯达克生物科技指数人民币C";           scale = "0.00"; pos = "90.20"; pct = "7.16"; mv = "0";      mvIsZero = $true;  rank = 1 },
    @{ idx = 4; code = "016471"; name = "广发纳斯达克生物科技指数美元C";            scale = "0.00"; pos = "90.20"; pct = "7.16"; mv = "0";      mvIsZero = $true;  rank = 1 }
)

foreach ($r in $rows) {
    $row = $r.idx + 2
    Set-NumberValue $newSheet.Cells.Item($row,1) $r.idx
    Set-TextValue   $newSheet.Cells.Item($row,2) $r.code
    Set-TextValue   $newSheet.Cells.Item($row,3) $r.name
    Set-TextValue   $newSheet.Cells.Item($row,4) $r.scale
    Set-TextValue   $newSheet.Cells.Item($row,5) $r.pos
    Set-TextValue   $newSheet.Cells.Item($row,6) $r.pct
    if ($r.mvIsZero) {
        Set-NumberValue $newSheet.Cells.Item($row,7) 0
    } else {
        Set-TextValue $newSheet.Cells.Item($row,7) $r.mv
    }
    Set-NumberValue $newSheet.Cells.Item($row,8) $r.rank
}

# Apply the same header/index-column style used by the neighbouring
# quarterly sheets (style "2" -> bordered/bold/centered header look).
$newSheet.Range("B1:H1").Style = $q3Sheet.Range("B1:H1").Style
$newSheet.Range("A2:A6").Style = $q3Sheet.Range("A2:A3").Style

# Re-select the sheet that was active/last before we started (the new
# sheet becomes selected automatically when added, which would otherwise
# move the persisted "tabSelected" flag away from the final sheet).
$wb.Worksheets.Item($lastSheetName).Activate()

# ---------------------------------------------------------------------
# 2. Insert a new row at the top of the "总计" table for the 2022-Q4
#    figures, pushing all the existing rows down by one.
# ---------------------------------------------------------------------
$summarySheet.Rows.Item(2).Insert()
$summarySheet.Range("A2:D2").ClearFormats()

Set-NumberValue $summarySheet.Cells.Item(2,1) 0
Set-TextValue   $summarySheet.Cells.Item(2,2) "2022-Q4"
Set-NumberValue $summarySheet.Cells.Item(2,3) 5
Set-NumberValue $summarySheet.Cells.Item(2,4) 0.32

$wb.Worksheets.Item($lastSheetName).Activate()
